$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title (D26)
$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

# Row 32: update title (D32) and link (E32)
$ws.Range("D32").Value = "이중차분법 (Difference In Difference)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/382"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "휴대폰 화면에 캡쳐 화면 넣기, 목업(mock-up) 활용 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%ED%9C%B4%EB%8C%80%ED%8F%B0-%ED%99%94%EB%A9%B4%EC%97%90-%EC%BA%A1%EC%B3%90-%ED%99%94%EB%A9%B4-%EB%84%A3%EA%B8%B0-%EB%AA%A9%EC%97%85mock-up-%ED%99%9C%EC%9A%A9-%EB%B0%A9%EB%B2%95"
